# The sheet contains a weekly price record table (rows 2-41). A new
# weekly observation is inserted as row 37 (pushing the former rows
# 37-41 down to 38-42), extending the used range to A1:T42.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 37, shifting existing rows 37-41 down to 38-42.
$ws.Rows.Item(37).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(37, 1).Value = 3
$ws.Cells.Item(37, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(37, 3).Value = "Coquimbo"
$ws.Cells.Item(37, 4).Value = 44476
$ws.Cells.Item(37, 5).Value = 5
$ws.Cells.Item(37, 6).Value = "Fruta"
$ws.Cells.Item(37, 7).Value = 100108
$ws.Cells.Item(37, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(37, 9).Value = 100108004
$ws.Cells.Item(37, 10).Value = "Papaya"
$ws.Cells.Item(37, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(37, 12).Value = "Primera"
$ws.Cells.Item(37, 13).Value = 68
$ws.Cells.Item(37, 14).Value = 17000
$ws.Cells.Item(37, 15).Value = 17000
$ws.Cells.Item(37, 16).Value = 17000
$ws.Cells.Item(37, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(37, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(37, 19).Value = 1700
$ws.Cells.Item(37, 20).Value = 10
